$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A21").Value = "Article review"
$ws.Range("B21").Value = "Nature Communications"
$ws.Range("C21").Value = 2025

$ws.Range("A22").Value = "Article review"
$ws.Range("B22").Value = "JSAN"
$ws.Range("C22").Value = "2023-2025"

$ws.Range("B28").Select()
